# Last minute changes to the LP functionality.
# Add a new sheet "Last 2 Days Calls" so the user can manually key in the
# on-call doctors for the last 2 days of the previous schedule month.

$wb = $excel.ActiveWorkbook
$roster = $wb.Worksheets.Item("Roster")
$publicHoliday = $wb.Worksheets.Item("Public Holiday")

# The Roster sheet ends up with a plain A1:B1 selection (and loses the
# "tabSelected"/scrolled-to-row-13 view it had before) once the new sheet
# becomes the active tab.
$roster.Activate()
$roster.Range("A1:B1").Select()

# Insert the new sheet after the last existing sheet ("Public Holiday").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Last 2 Days Calls"

# ---- Header row: reuse the look of the Roster header, then set the text ----
$roster.Range("A1:C1").Copy()
$newSheet.Range("A1:C1").PasteSpecial(-4122)
$newSheet.Range("A1").Value = "Email"
$newSheet.Range("B1").Value = "Name"
$newSheet.Range("C1").Value = "Duty Date"

# ---- Data rows: Email / Name columns reuse the Roster's bordered style ----
$roster.Range("A2:B2").Copy()
$newSheet.Range("A2:B6").PasteSpecial(-4122)

# Duty Date column reuses a cell that is already formatted as a date.
$publicHoliday.Range("A2").Copy()
$newSheet.Range("C2:C6").PasteSpecial(-4122)

$newSheet.Range("A2").Value = "V"
$newSheet.Range("B2").Value = "V"
$newSheet.Range("C2").Value = "2/28/2020"

$newSheet.Range("A3").Value = "A"
$newSheet.Range("B3").Value = "A"
$newSheet.Range("C3").Value = "2/28/2020"

$newSheet.Range("A4").Value = "B"
$newSheet.Range("B4").Value = "B"
$newSheet.Range("C4").Value = "2/29/2020"

$newSheet.Range("A5").Value = "W"
$newSheet.Range("B5").Value = "W"
$newSheet.Range("C5").Value = "2/29/2020"

$newSheet.Range("A6").Value = "X"
$newSheet.Range("B6").Value = "X"
$newSheet.Range("C6").Value = "2/29/2020"

$newSheet.Columns("C:C").AutoFit()

# Leave the new sheet active/selected, matching the tab the author was last on.
$newSheet.Activate()
$newSheet.Range("G26").Select()
